# Apply the "performances_polar_apolar" edit:
#  - sheets "apolar anisotropy" / "apolar isotropy": headers change, row-labels
#    (svr/poly/xgb/nn/stacked) added in column A, values shifted to B:E.
#  - sheets "polar anisotropy" / "polar isotropy": same header text kept, but
#    headers shift from A:D to B:E, with row-labels (poly/svr/xgb/nn/stacked)
#    added in column A and values shifted right.
#  - two new sheets "all anisotropy" / "all isotropy" appended at the end,
#    built the same way as the polar sheets.

$wb = $excel.ActiveWorkbook

# A cell that already carries the bold/border/center header style, used as a
# style-only donor for every header / row-label cell we touch below (copying
# from a cell on the *same* still-being-built sheet would be a no-op on a
# brand-new blank sheet, so we always pull from this pre-existing sheet).
$styleDonorSheet = $wb.Worksheets.Item("apolar isotropy")
$styleDonorCell = $styleDonorSheet.Range("B1")

function Set-SheetPageProps($ws) {
    # Match the outline / page-margin conventions already used by the
    # existing sheets (summaryBelow/summaryRight outline + 0.75/1/0.5in margins).
    $ws.Outline.SummaryRow = 1
    $ws.Outline.SummaryColumn = 1
    $ws.PageSetup.LeftMargin = 54
    $ws.PageSetup.RightMargin = 54
    $ws.PageSetup.TopMargin = 72
    $ws.PageSetup.BottomMargin = 72
    $ws.PageSetup.HeaderMargin = 36
    $ws.PageSetup.FooterMargin = 36
}

function Fill-Sheet($ws, $headers, $rowLabels, $data) {
    # Make sure the cell that used to hold the single A1 header text is fully
    # cleared (style included) before anything else -- its old value moves to B1.
    $ws.Range("A1").Clear()

    $cols = @("B", "C", "D", "E")

    for ($i = 0; $i -lt $headers.Count; $i++) {
        $cell = $ws.Range("$($cols[$i])1")
        # Copy format (bold/border/center) from an already-styled header cell,
        # then overwrite with the new text -- Copy preserves the existing
        # style index instead of minting a near-duplicate one.
        $styleDonorCell.Copy($cell)
        $cell.Value = $headers[$i]
    }

    for ($r = 0; $r -lt $rowLabels.Count; $r++) {
        $rowNum = $r + 2
        $labelCell = $ws.Range("A$rowNum")
        $styleDonorCell.Copy($labelCell)
        $labelCell.Value = $rowLabels[$r]

        $rowValues = $data[$r]
        for ($c = 0; $c -lt $rowValues.Count; $c++) {
            $ws.Range("$($cols[$c])$rowNum").Value = $rowValues[$c]
        }
    }
}

# ---------------------------------------------------------------------------
# apolar anisotropy
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("apolar anisotropy")
Fill-Sheet $ws `
    @("Ei + Alpha + Pi", "Ei + Alpha", "Pi + Ei", "Pi + Alpha") `
    @("svr", "poly", "xgb", "nn", "stacked") `
    @(
        @(0.5265497639688039, 0.9610753334149637, 1.271563834511352, 0.4582919031113932),
        @(0.7539364898971267, 0.7465780862731102, 0.9783564291005794, 0.6453247856342541),
        @(0.8482818073022681, 1.082946839015822, 0.9111777761615172, 0.5627059711741039),
        @(0.3826688805714364, 1.170002279825135, 1.254272665982813, 0.4502129087680733),
        @(0.51639674670516, 0.9637088664367368, 1.049129694500607, 0.4462339865424551)
    )

# ---------------------------------------------------------------------------
# apolar isotropy
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("apolar isotropy")
Fill-Sheet $ws `
    @("Ei + Alpha + Pi", "Ei + Alpha", "Pi + Ei", "Pi + Alpha") `
    @("svr", "poly", "xgb", "nn", "stacked") `
    @(
        @(0.5626601962334903, 0.932495725147756, 1.271563834511352, 0.4636680943493121),
        @(0.5271616330798463, 1.011789946258571, 0.9783564291005794, 0.3870707482702535),
        @(0.9426680533400128, 0.9694174411579096, 0.9111777761615172, 0.4899423608629933),
        @(0.3190289668897709, 0.8486709838753951, 1.254272665982813, 0.4234874674772652),
        @(0.488580156572951, 0.9063687251534668, 1.049129694500607, 0.3882818462850995)
    )

# ---------------------------------------------------------------------------
# polar anisotropy
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("polar anisotropy")
Fill-Sheet $ws `
    @("Ei + Alpha + Dipole + Pi", "Ei + Alpha + Dipole", "Alpha + Dipole", "Alpha + Dipole + Pi") `
    @("poly", "svr", "xgb", "nn", "stacked") `
    @(
        @(0.2108254817370543, 0.3793491622610133, 0.3418258538616912, 0.2154622538444675),
        @(0.3055489869880953, 0.3320233910124897, 0.5175891612008552, 0.3099075827789245),
        @(0.2534800511804198, 0.2809519345643381, 0.375393804398327, 0.3891850978305191),
        @(0.2541444994724172, 0.3690626397914039, 0.3820350294391228, 0.2892566927193675),
        @(0.2076097802159532, 0.2524215509446512, 0.3620969116592002, 0.2657215067766542)
    )

# ---------------------------------------------------------------------------
# polar isotropy
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("polar isotropy")
Fill-Sheet $ws `
    @("Ei + Alpha + Dipole + Pi", "Ei + Alpha + Dipole", "Alpha + Dipole", "Alpha + Dipole + Pi") `
    @("poly", "svr", "xgb", "nn", "stacked") `
    @(
        @(0.3321349716619197, 0.4477577797951195, 0.4083421828570664, 0.4408974775179098),
        @(0.3638766402041535, 0.4118322651507619, 0.4351820676343802, 0.4140239721744264),
        @(0.3210148473180056, 0.3324928769631271, 0.3671071970855833, 0.3585627950437669),
        @(0.3440934470386112, 0.3627337999359304, 0.3611304984549407, 0.3683702957459811),
        @(0.2877024134822573, 0.3325173805703944, 0.3590381502466007, 0.3595235085012145)
    )

# ---------------------------------------------------------------------------
# new sheet: all anisotropy (appended after polar isotropy)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsAllAniso = $wb.Worksheets.Add($null, $lastSheet)
$wsAllAniso.Name = "all anisotropy"
Set-SheetPageProps $wsAllAniso
Fill-Sheet $wsAllAniso `
    @("Ei + Alpha + Dipole + Pi", "Ei + Alpha + Dipole", "Alpha + Dipole", "Alpha + Dipole + Pi") `
    @("poly", "svr", "xgb", "nn", "stacked") `
    @(
        @(0.2152276859145625, 0.59951976322566, 0.4020997742751305, 0.1948498644656363),
        @(0.2637070494491611, 0.8449949146977472, 0.4947970364806629, 0.3224539525264279),
        @(0.3707633288011138, 0.4662531522940735, 0.457587715842526, 0.4289947940047193),
        @(0.2651304380493411, 0.5395860301353854, 0.4866719311209929, 0.2775617324526775),
        @(0.2000430580789068, 0.5689972952765208, 0.4197662330536388, 0.2520192417944014)
    )

# ---------------------------------------------------------------------------
# new sheet: all isotropy (appended after all anisotropy)
# ---------------------------------------------------------------------------
$wsAllIso = $wb.Worksheets.Add($null, $wsAllAniso)
$wsAllIso.Name = "all isotropy"
Set-SheetPageProps $wsAllIso
Fill-Sheet $wsAllIso `
    @("Ei + Alpha + Dipole + Pi", "Ei + Alpha + Dipole", "Alpha + Dipole", "Alpha + Dipole + Pi") `
    @("poly", "svr", "xgb", "nn", "stacked") `
    @(
        @(0.3470499220251377, 0.5040248465377545, 0.3650544642875731, 0.3441594007058143),
        @(0.3504179712992437, 0.5823934342922048, 0.4141216693794856, 0.3166355566482289),
        @(0.3106855463025898, 0.4769994961412236, 0.4384711290856459, 0.428689951218051),
        @(0.5304828493332466, 0.4778160486610435, 0.4053688631875019, 0.3869492856871732),
        @(0.3395245737205701, 0.4451432876186266, 0.3609772571659764, 0.3256863613147857)
    )

Write-Output "edit complete"
